$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh the export timestamps on the "data" sheet (F2:F26) ---
$times = @(
  "2021-10-05 14:20:11.659131",
  "2021-10-05 14:20:11.659140",
  "2021-10-05 14:20:11.659143",
  "2021-10-05 14:20:11.659146",
  "2021-10-05 14:20:11.659151",
  "2021-10-05 14:20:11.659153",
  "2021-10-05 14:20:11.659156",
  "2021-10-05 14:20:11.659159",
  "2021-10-05 14:20:11.659162",
  "2021-10-05 14:20:11.659164",
  "2021-10-05 14:20:11.659167",
  "2021-10-05 14:20:11.659170",
  "2021-10-05 14:20:11.659172",
  "2021-10-05 14:20:11.659175",
  "2021-10-05 14:20:11.659178",
  "2021-10-05 14:20:11.659180",
  "2021-10-05 14:20:11.659184",
  "2021-10-05 14:20:11.659187",
  "2021-10-05 14:20:11.659189",
  "2021-10-05 14:20:11.659192",
  "2021-10-05 14:20:11.659195",
  "2021-10-05 14:20:11.659197",
  "2021-10-05 14:20:11.659200",
  "2021-10-05 14:20:11.659203",
  "2021-10-05 14:20:11.659206"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $times[$i]
}

# --- 2. Add the new "metadata" worksheet after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Copy the header formatting (bold, border, centered) from the data sheet
# so the new header row reuses the same cell style.
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$dataSheet.Range("F1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Copy the style for the numeric index cell A2 (re-use style "A2" from data sheet).
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Header row ---
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# --- 4. Data row ---
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 2).Value = "Familial cicatricial alopecia"
$metaSheet.Cells.Item(2, 3).Value = 305
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "1.3"
$metaSheet.Cells.Item(2, 4).Style = "Normal"
$metaSheet.Cells.Item(2, 5).Value = "2021-08-31T11:16:38.259932Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:20:11.655348"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/305/?format=json"

$metaSheet.Range("A1").Select() | Out-Null
